$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Contacts"

# Clear old content
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "quantity"

$data = @(
    "E23-11-Premium",
    "E25-11-Premium",
    "E27-11.2-Premium",
    "E28-11.2-Premium",
    "E29-11.2-Premium",
    "E08-10.4-Premium"
)

$r = 2
foreach ($name in $data) {
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = 1
    $r = $r + 1
}

$ws.Range("C1").Value = "serials"

# Header formatting: dark blue fill, white Arial font
$header = $ws.Range("A1:C1")
$header.Font.Name = "Arial"
$header.Font.Family = 1
$header.Font.Color = 16777215
$header.Interior.Color = 9851952

Write-Host "data written"
